# Tracker update for 13 Feb 2022
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 8 with the new data point (09-Feb-22)
$ws.Range("A8").Value = 44601
$ws.Range("B8").Value = 0.04
$ws.Range("C8").Value = 0.14000000000000001
$ws.Range("D8").Value = 0.05
$ws.Range("E8").Value = 0.52
$ws.Range("F8").Value = 0.08

# Update the view to reflect the new scroll/selection position
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("E7").Select()
